$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right count and Wrong penalty
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total, Wrong total, and Max display text
$ws.Range("B12").Value = 36
$ws.Range("C12").Value = -16
$ws.Range("E12").Value = "20 / 112"
